# UserStories.xlsx - "Tasks" sheet: mark the "Klasse Login erstellen" task
# (row 18) as done, recording its completion date - mirroring the other
# finished tasks above it (e.g. row 15).
#
# Once C18 no longer holds the placeholder status abbreviation "b", that
# shared string has no remaining references anywhere in the workbook, so
# Excel drops it from sharedStrings.xml on save (all following shared
# string indices shift down by one as a side effect - rows 19-22 keep
# their original text, just pointing at the renumbered entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Row 18: "Klasse Login erstellen" task is now finished -> status "done"
# with a completion date (2019-03-01, same date serial as row 15).
$ws.Range("C18").Value = "done"
$ws.Range("D18").Value2 = 43525

# Copy the date-cell formatting (wrapped text, short-date number format)
# from the existing "done" row (D15) onto the newly populated D18 cell.
$ws.Range("D15").Copy()
$ws.Range("D18").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection to the newly edited cell, like a user would
# after typing the completion date.
$ws.Range("D18").Select()
